$d = $word.ActiveDocument

# --- 1. "Intor carrer .tops site." paragraph: merge the proofErr-fragmented
#        runs ("Intor" / " " / "carrer" / " .tops site.") into a single run,
#        dropping the spell-check proofErr markers entirely.
$pIntor = $d.Paragraphs.Item(3)
$rIntor = $pIntor.Range.Duplicate
[void]$rIntor.MoveStart(1, -1)
[void]$rIntor.Find.Execute("`rIntor carrer .tops site.", $false, $false, $false, $false, $false, $true, 1, $false, "`rIntor carrer .tops site.", 2)

# --- 2. ".edu" paragraph: merge "." + proofErr("edu") into a single run.
$pEdu = $d.Paragraphs.Item(11)
$rEdu = $pEdu.Range.Duplicate
[void]$rEdu.MoveStart(1, -1)
[void]$rEdu.Find.Execute("`r.edu", $false, $false, $false, $false, $false, $true, 1, $false, "`r.edu", 2)

# --- 3. Add the new "float type c code" entries after "Git hub - git bash".
$pGit = $d.Paragraphs.Item(17)
$rGit = $pGit.Range.Duplicate
[void]$rGit.InsertAfter("`r19-01-26`r`rPrograming language `rDev c++-code soirce")

# The blank line between "19-01-26" and "Programing language " should stay a
# truly empty paragraph (no run), matching the rest of the document's blank
# paragraphs.
$pBlank = $d.Paragraphs.Item(19)
[void]$pBlank.Range.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'/>")
